$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same style used by the header row (A9) to the row-label
# cells in A10, A11, A12 -- copy formats only so the existing style
# index is reused instead of minting a new one.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the marking/total values.
# C11 stores its value as text (e.g. "-3"), so force it back in as text
# (leading apostrophe) rather than letting it auto-convert to a number.
$ws.Range("C11").Formula = "'-1"
$ws.Range("C12").Value = -13
$ws.Range("E12").Value = "27/140"
